$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the three new ones ---------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "indexedList"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "indexedListAsLeaf"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "indexedListAsLeafTestOption"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "indexedListAsLeafTestLeft"

# --- Sheet2: indexedListAsLeaf ------------------------------------------
$ws2.Range("A1").Value = "####"
$ws2.Range("C1").Value = "listAsLeaf#test"
$ws2.Range("D1").Value = "listAsLeaf#list[0]"
$ws2.Range("E1").Value = "listAsLeaf#list[1]"
$ws2.Range("F1").Value = "listAsLeaf#list[2]"

$ws2.Range("A3").Value = "listAsLeaf#~"
$ws2.Range("C3").Value = "aaa"
$ws2.Range("D3").Value = "bbb"
$ws2.Range("E3").Value = "ccc"
$ws2.Range("F3").Value = "ddd"

$ws2.Range("D4").Value = 1
$ws2.Range("E4").Value = 2
$ws2.Range("F4").Value = 3

$ws2.Range("D5").Value = 4
$ws2.Range("F5").Value = 5

$ws2.Range("E6").Value = 6

$ws2.Columns.Item(3).ColumnWidth = 13.875
$ws2.Range("D1:F1").Columns.ColumnWidth = 15.6875

# --- Sheet3: indexedListAsLeafTestOption --------------------------------
$ws3.Range("A1").Value = "####"
$ws3.Range("C1").Value = "listAsLeafTestOption#test?type=string"
$ws3.Range("D1").Value = "listAsLeafTestOption#list[0]"
$ws3.Range("E1").Value = "listAsLeafTestOption#list[1]?type=string"
$ws3.Range("F1").Value = "listAsLeafTestOption#list[2]"

$ws3.Range("A3").Value = "listAsLeafTestOption#~"
$ws3.Range("C3").Value = 1
$ws3.Range("D3").Value = "bbb"
$ws3.Range("E3").Value = "ccc"
$ws3.Range("F3").Value = "ddd"

$ws3.Range("D4").Value = 1
$ws3.Range("E4").Value = 2
$ws3.Range("F4").Value = 3

$ws3.Range("D5").Value = 4
$ws3.Range("F5").Value = 5

$ws3.Range("E6").Value = 6

$ws3.Columns.Item(1).ColumnWidth = 21.625
$ws3.Columns.Item(2).ColumnWidth = 2.375
$ws3.Columns.Item(3).ColumnWidth = 34.8125
$ws3.Columns.Item(4).ColumnWidth = 25.4375
$ws3.Columns.Item(5).ColumnWidth = 36.625
$ws3.Columns.Item(6).ColumnWidth = 25.4375

# --- Sheet4: indexedListAsLeafTestLeft ----------------------------------
$ws4.Range("A1").Value = "####"
$ws4.Range("C1").Value = "listAsLeafTestLeft[0]#test?type=string"
$ws4.Range("D1").Value = "listAsLeafTestLeft[0]#list[0]"
$ws4.Range("E1").Value = "listAsLeafTestLeft[0]#list[1]?type=string"
$ws4.Range("F1").Value = "listAsLeafTestLeft[0]#list[2]"
$ws4.Range("H1").Value = "listAsLeafTestLeft[1]#test"
$ws4.Range("I1").Value = "listAsLeafTestLeft[1]#list[0]"
$ws4.Range("J1").Value = "listAsLeafTestLeft[1]#list[1]?type=string"
$ws4.Range("K1").Value = "listAsLeafTestLeft[1]#list[2]"

$ws4.Range("A3").Value = "listAsLeafTestLeft[0]#~, listAsLeafTestLeft[1]#~"
$ws4.Range("C3").Value = 1
$ws4.Range("D3").Value = "bbb"
$ws4.Range("E3").Value = "ccc"
$ws4.Range("F3").Value = "ddd"
$ws4.Range("H3").Value = 2
$ws4.Range("I3").Value = 3
$ws4.Range("J3").Value = 4
$ws4.Range("K3").Value = 5

$ws4.Range("D4").Value = 1
$ws4.Range("E4").Value = 2
$ws4.Range("F4").Value = 3

$ws4.Range("D5").Value = 4
$ws4.Range("F5").Value = 5

$ws4.Range("E6").Value = 6

$ws4.Columns.Item(1).ColumnWidth = 21.625
$ws4.Columns.Item(2).ColumnWidth = 2.375
$ws4.Columns.Item(3).ColumnWidth = 34.8125
$ws4.Columns.Item(4).ColumnWidth = 25.4375
$ws4.Columns.Item(5).ColumnWidth = 36.625
$ws4.Columns.Item(6).ColumnWidth = 25.4375
$ws4.Columns.Item(7).ColumnWidth = 6.1875
$ws4.Columns.Item(8).ColumnWidth = 23.75
$ws4.Columns.Item(9).ColumnWidth = 25.5625
$ws4.Columns.Item(10).ColumnWidth = 36.8125
$ws4.Columns.Item(11).ColumnWidth = 25.5625

# --- Selections / active sheet ------------------------------------------
$ws1.Activate()
$ws1.Range("C11").Select() | Out-Null

$ws2.Activate()
$ws2.Range("E2").Select() | Out-Null

$ws3.Activate()
$ws3.Range("C3").Select() | Out-Null

$ws4.Activate()
$ws4.Range("H3").Select() | Out-Null
